# Scheduled-runner refresh of per-item market/profit figures on each
# crafting-job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) in Bahamut_Profits.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) for the leves whose market data changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 323609.5
$ws.Range("I19").Value = 387978.06
$ws.Range("J19").Value = 1766.6666
$ws.Range("K19").Value = 387978.06
$ws.Range("L19").Value = 1766.6666
$ws.Range("M19").Value = -387803.06
$ws.Range("N19").Value = -2116.6666

$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -5126
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -13632
$ws.Range("N72").ClearContents()

$ws.Range("H125").Value = 1265.6666
$ws.Range("I125").Value = 1404.5714
$ws.Range("J125").Value = 1071.2
$ws.Range("K125").Value = 12641.1426
$ws.Range("L125").Value = 9640.800000000001
$ws.Range("M125").Value = -10181.1426
$ws.Range("N125").Value = -14560.8

$ws.Range("H137").Value = 949.25
$ws.Range("I137").Value = 915.6667
$ws.Range("J137").Value = 1050
$ws.Range("K137").Value = 2747.0001
$ws.Range("L137").Value = 3150
$ws.Range("M137").Value = -197.0001000000002
$ws.Range("N137").Value = -8250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 10514.286
$ws.Range("I31").Value = 3966.6667
$ws.Range("J31").Value = 49800
$ws.Range("K31").Value = 3966.6667
$ws.Range("L31").Value = 49800
$ws.Range("M31").Value = -3672.6667
$ws.Range("N31").Value = -50388

$ws.Range("H32").Value = 6639.2437
$ws.Range("I32").Value = 5365.3384
$ws.Range("J32").Value = 13008.77
$ws.Range("K32").Value = 5365.3384
$ws.Range("L32").Value = 13008.77
$ws.Range("M32").Value = -5078.3384
$ws.Range("N32").Value = -13582.77

$ws.Range("H102").Value = 3408.889
$ws.Range("I102").Value = 3525
$ws.Range("J102").Value = 2480
$ws.Range("K102").Value = 3525
$ws.Range("L102").Value = 2480
$ws.Range("M102").Value = -1903
$ws.Range("N102").Value = -5724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 17906
$ws.Range("J81").Value = 17906
$ws.Range("L81").Value = 17906
$ws.Range("N81").Value = -20028

$ws.Range("H84").Value = 17906
$ws.Range("J84").Value = 17906
$ws.Range("L84").Value = 53718
$ws.Range("N84").Value = -64326

$ws.Range("H107").Value = 10554.8125
$ws.Range("I107").Value = 1861.5454
$ws.Range("J107").Value = 29680
$ws.Range("K107").Value = 1861.5454
$ws.Range("L107").Value = 29680
$ws.Range("M107").Value = 58.45460000000003
$ws.Range("N107").Value = -33520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 16339.667
$ws.Range("I26").Value = 6019
$ws.Range("K26").Value = 6019
$ws.Range("M26").Value = -5732

$ws.Range("H132").Value = 1538.5
$ws.Range("I132").Value = 1103
$ws.Range("J132").Value = 3498.25
$ws.Range("K132").Value = 3309
$ws.Range("L132").Value = 10494.75
$ws.Range("M132").Value = -779
$ws.Range("N132").Value = -15554.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1280.6
$ws.Range("I114").Value = 875.75
$ws.Range("J114").Value = 2900
$ws.Range("K114").Value = 2627.25
$ws.Range("L114").Value = 8700
$ws.Range("M114").Value = 626.75
$ws.Range("N114").Value = -15208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4990.4614
$ws.Range("I70").Value = 5201.1665
$ws.Range("J70").Value = 4809.857
$ws.Range("K70").Value = 5201.1665
$ws.Range("L70").Value = 4809.857
$ws.Range("M70").Value = -4931.1665
$ws.Range("N70").Value = -5349.857

$ws.Range("H73").Value = 4990.4614
$ws.Range("I73").Value = 5201.1665
$ws.Range("J73").Value = 4809.857
$ws.Range("K73").Value = 5201.1665
$ws.Range("L73").Value = 4809.857
$ws.Range("M73").Value = -4265.1665
$ws.Range("N73").Value = -6681.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2463.625
$ws.Range("I7").Value = 2617.3333
$ws.Range("J7").Value = 2002.5
$ws.Range("K7").Value = 2617.3333
$ws.Range("L7").Value = 2002.5
$ws.Range("M7").Value = -2505.3333
$ws.Range("N7").Value = -2226.5

$ws.Range("H40").Value = 2527302.2
$ws.Range("I40").Value = 10101010
$ws.Range("J40").Value = 2733
$ws.Range("K40").Value = 10101010
$ws.Range("L40").Value = 2733
$ws.Range("M40").Value = -10100874
$ws.Range("N40").Value = -3005

$ws.Range("H100").Value = 1754.125
$ws.Range("I100").Value = 1713.8334
$ws.Range("J100").Value = 1875
$ws.Range("K100").Value = 1713.8334
$ws.Range("L100").Value = 1875
$ws.Range("M100").Value = -1172.8334
$ws.Range("N100").Value = -2957

$ws.Range("H126").Value = 2463.625
$ws.Range("I126").Value = 2617.3333
$ws.Range("J126").Value = 2002.5
$ws.Range("K126").Value = 7851.999899999999
$ws.Range("L126").Value = 6007.5
$ws.Range("M126").Value = -5381.999899999999
$ws.Range("N126").Value = -10947.5

$ws.Range("H132").Value = 4391.5
$ws.Range("I132").Value = 6566.6665
$ws.Range("J132").Value = 3666.4443
$ws.Range("K132").Value = 19699.9995
$ws.Range("L132").Value = 10999.3329
$ws.Range("M132").Value = -17169.9995
$ws.Range("N132").Value = -16059.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1072.1428
$ws.Range("I100").Value = 1242.7273
$ws.Range("J100").Value = 446.66666
$ws.Range("K100").Value = 2485.4546
$ws.Range("L100").Value = 893.33332
$ws.Range("M100").Value = -1944.4546
$ws.Range("N100").Value = -1975.33332

$ws.Range("H126").Value = 904.0714
$ws.Range("I126").Value = 735.2
$ws.Range("J126").Value = 1326.25
$ws.Range("K126").Value = 2205.6
$ws.Range("L126").Value = 3978.75
$ws.Range("M126").Value = 264.3999999999996
$ws.Range("N126").Value = -8918.75

$ws.Range("H132").Value = 849.8823
$ws.Range("I132").Value = 705.931
$ws.Range("J132").Value = 1684.8
$ws.Range("K132").Value = 2117.793
$ws.Range("L132").Value = 5054.4
$ws.Range("M132").Value = 412.2069999999999
$ws.Range("N132").Value = -10114.4

Write-Output "Applied all Bahamut_Profits updates"